$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-10 (columns B-G)
$data = @{
    2  = @(0.2240261301401967, 1.34776042189687, 8.330079577401087, 2.886187723867089, 2.9061124733609, 51)
    3  = @(0.3015350435189421, 1.634473129698246, 9.405841710792368, 3.066894473370802, 3.083021035019027, 50)
    4  = @(0.2397024084672859, 1.273559806460503, 6.25558203749066, 2.501116158336246, 2.515403027889767, 49)
    5  = @(0.28645024546988, 1.336241691786354, 6.93067963156568, 2.632618398394587, 2.644681707289316, 48)
    6  = @(0.2510080378810327, 1.618237384748787, 9.429041488443973, 3.070674435436615, 3.09348443368957, 47)
    7  = @(0.317070059663221, 1.568161529071841, 8.05755616701061, 2.838583478957525, 2.851989739099063, 46)
    8  = @(0.2722875540392151, 1.725110569588582, 8.603484540338279, 2.933169708751657, 2.953505164405249, 45)
    9  = @(0.3499353129240411, 1.708867493651148, 8.79813827900456, 2.966165585230292, 2.979503924895169, 44)
    10 = @(0.2573786280145778, 1.691059197863888, 8.352307215966194, 2.890035850290822, 2.912619213524731, 43)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

# New row 11: Q9 label plus its stats
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = "Q9"

$ws.Cells.Item(11, 2).Value = 0.2874384264431032
$ws.Cells.Item(11, 3).Value = 1.714766158501543
$ws.Cells.Item(11, 4).Value = 8.919363088342259
$ws.Cells.Item(11, 5).Value = 2.986530275812093
$ws.Cells.Item(11, 6).Value = 3.008699476163549
$ws.Cells.Item(11, 7).Value = 42
